$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.382.50"
$ws.Range("D3").Value = "3.674.76"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'645.08"
$ws.Range("E5").Value = "  -5.22%  "
$ws.Range("D6").Value = "'159.91"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.498"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "'7.09"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").Value = "'0.450"
$ws.Range("D12").Value = "'0.0000232"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "4.291.81"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "'32.75"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "3.645.31"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "69.358.61"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D18").Value = "'16.02"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'466.57"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "'9.94"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "'0.647"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "'79.47"
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "3.820.95"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'9.06"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  -2.85%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "'26.90"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "3.667.32"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'8.43"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("E39").Value = "  -5.70%  "
$ws.Range("D40").Value = "'178.73"
$ws.Range("E40").Value = "  +4.42%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "'0.0901"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "'0.926"
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'27.36"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "'0.000271"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  -4.30%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.85"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'1.25"
$ws.Range("E51").Value = "  -3.85%  "
